$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update phone numbers in column B (B2:B7) with the new value
$ws.Range("B2:B7").Value = "5512990000000"

# Update selection to match the new active range
$ws.Range("B2:B7").Select()
